# Apply version bump, clear "Experimental" value, and update Date on the
# "Metadata" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version: 1.0.4 -> 1.0.7
$ws.Range("B3").Value = "1.0.7"

# Experimental value "false" is removed entirely (cell becomes empty)
$ws.Range("B7").ClearContents()

# Date: 2025-04-11 -> 2025-09-12
# Force text so Excel doesn't auto-convert the ISO-like date string into a
# real date serial number (it must stay a literal text value).
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "2025-09-12"
